$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at the top, shifting existing rows down
$ws.Rows.Item(1).Insert()

# Populate the new header cells
$ws.Range("A1").Value = "xxxx"
$ws.Range("B1").Value = "yyyy"

# Restore selection like the target file (B1 selected)
$ws.Range("B1").Select()
